# Lectura y Validación de Movimientos del XML generando Objetos tipo Movimiento
#
# The first data row (row 2) held placeholder/test data ("111111111111" /
# "Prueba 1") in columns A-C (NroContrato, CUIT, DENOMINACION) that no
# longer corresponds to a real parsed Movimiento. Clear those three cells
# while leaving the rest of the row (DOMICILIO, CODIGOPOSTAL, PRODUCTOR)
# and every other row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:C2").ClearContents()
